# Add a "date_of_birth" column (column E) to Sheet1 holding each patient's
# birth date, formatted as a short date - mirroring the existing
# patient_id/wbc_value columns added earlier in the workbook. Leaves the
# selection where the user's cursor ended up after entering the data (one
# row below the table).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Birth-date serials (Excel 1900 date system) for the four patients, in the
# same row order as the existing A2:A5 patient_id values.
# 16755 -> 11/14/1945, 33197 -> 11/20/1990, 28247 -> 5/2/1977, 38756 -> 2/8/2006
$dobSerials = @(16755, 33197, 28247, 38756)

# Write + format the first data cell, then propagate that cell's formatting
# (short-date number format) down the column via copy/paste so every data
# cell shares a single style entry - the same as Excel does when you format
# one cell and fill/paste it down a column.
$ws.Range("E2").Value = $dobSerials[0]
$ws.Range("E2").NumberFormat = "mm-dd-yy"
$ws.Range("E2").Copy()
$ws.Range("E3:E5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("E3").Value = $dobSerials[1]
$ws.Range("E4").Value = $dobSerials[2]
$ws.Range("E5").Value = $dobSerials[3]

# Auto-size the new column to fit its (date) contents.
$ws.Columns.Item(5).AutoFit()

# Header goes in last, as plain text like "patient_id"/"wbc_value".
$ws.Range("E1").Value = "date_of_birth"

# Cursor ends up just below the table after the edits.
[void]$ws.Range("E6").Select()
